$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Personas legisladoras: mark existing rows as "Completa" (status column E) ---
$ws.Range("E20:E30").Value = "Completa"

# --- New review comments (column F) ---
$ws.Range("F19").Value = "Esto se debe a los procesos internos del aplicativo, se optimizará al final del proyecto."
$ws.Range("F31:F32").Value = "falta guardado"

# --- E18 loses its top border (cosmetic cleanup) ---
$ws.Range("E18").Borders.Item(8).LineStyle = -4142

# --- Extend the "Completa/Pendiente" conditional formatting down through row 57 ---
$cf = $ws.Range("E3:E17").FormatConditions
$cf.Item(1).ModifyAppliesToRange($ws.Range("E3:E57"))

# --- Extend the list data-validation on column E down through row 57 ---
$ws.Range("E18:E57").Validation.Add(3, 1, 1, "Pendiente, Completa")

# --- Update the selection / scroll position to reflect work down at row 31 ---
$ws.Range("F31").Select()
